$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 5055830.5
$ws.Cells.Item(106, 9).Value = 6498930
$ws.Cells.Item(106, 11).Value = 6498930
$ws.Cells.Item(106, 13).Value = -6498299
$ws.Cells.Item(112, 8).Value = 9100.9375
$ws.Cells.Item(112, 10).Value = 9100.9375
$ws.Cells.Item(112, 12).Value = 27302.8125
$ws.Cells.Item(112, 14).Value = -29518.8125
$ws.Cells.Item(134, 8).Value = 29852.133
$ws.Cells.Item(134, 10).Value = 29852.133
$ws.Cells.Item(134, 12).Value = 29852.133
$ws.Cells.Item(134, 14).Value = -39992.133
$ws.Cells.Item(138, 8).Value = 3455.5615
$ws.Cells.Item(138, 10).Value = 3104.2163
$ws.Cells.Item(138, 12).Value = 9312.6489
$ws.Cells.Item(138, 14).Value = -19592.6489

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2470.691
$ws.Cells.Item(32, 9).Value = 2146.0925
$ws.Cells.Item(32, 10).Value = 19999
$ws.Cells.Item(32, 11).Value = 2146.0925
$ws.Cells.Item(32, 12).Value = 19999
$ws.Cells.Item(32, 13).Value = -1859.0925
$ws.Cells.Item(32, 14).Value = -20573
$ws.Cells.Item(45, 8).Value = 2684.4614
$ws.Cells.Item(45, 9).Value = 2155.2856
$ws.Cells.Item(45, 11).Value = 2155.2856
$ws.Cells.Item(45, 13).Value = -1778.2856
$ws.Cells.Item(102, 8).Value = 18727
$ws.Cells.Item(102, 9).Value = 18727
$ws.Cells.Item(102, 11).Value = 18727
$ws.Cells.Item(102, 13).Value = -17105
$ws.Cells.Item(132, 8).Value = 2729.48
$ws.Cells.Item(132, 9).Value = 2709.0217
$ws.Cells.Item(132, 11).Value = 8127.0651
$ws.Cells.Item(132, 13).Value = -5597.0651

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 22583.555
$ws.Cells.Item(94, 9).Value = 530.6
$ws.Cells.Item(94, 10).Value = 50149.75
$ws.Cells.Item(94, 11).Value = 530.6
$ws.Cells.Item(94, 12).Value = 50149.75
$ws.Cells.Item(94, 13).Value = -79.60000000000002
$ws.Cells.Item(94, 14).Value = -51051.75
$ws.Cells.Item(99, 8).Value = 1811
$ws.Cells.Item(99, 9).Value = 1767.9375
$ws.Cells.Item(99, 11).Value = 1767.9375
$ws.Cells.Item(99, 13).Value = -269.9375

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 1750600
$ws.Cells.Item(6, 9).Value = 1750600
$ws.Cells.Item(6, 11).Value = 1750600
$ws.Cells.Item(6, 13).Value = -1750487
$ws.Cells.Item(25, 8).Value = 3000
$ws.Cells.Item(25, 9).Value = 3000
$ws.Cells.Item(25, 11).Value = 3000
$ws.Cells.Item(25, 13).Value = -2826
$ws.Cells.Item(86, 8).Value = 7044.8887
$ws.Cells.Item(86, 9).Value = 5580.6
$ws.Cells.Item(86, 11).Value = 5580.6
$ws.Cells.Item(86, 13).Value = -4457.6
$ws.Cells.Item(89, 8).Value = 7044.8887
$ws.Cells.Item(89, 9).Value = 5580.6
$ws.Cells.Item(89, 11).Value = 27903
$ws.Cells.Item(89, 13).Value = -22287
$ws.Cells.Item(132, 8).Value = 4936
$ws.Cells.Item(132, 9).Value = 4472.125
$ws.Cells.Item(132, 10).Value = 6791.5
$ws.Cells.Item(132, 11).Value = 13416.375
$ws.Cells.Item(132, 12).Value = 20374.5
$ws.Cells.Item(132, 13).Value = -10886.375
$ws.Cells.Item(132, 14).Value = -25434.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 129244.46
$ws.Cells.Item(37, 10).Value = 129244.46
$ws.Cells.Item(37, 12).Value = 387733.38
$ws.Cells.Item(37, 14).Value = -387957.38
$ws.Cells.Item(64, 8).Value = 1001756.3
$ws.Cells.Item(64, 10).Value = 1168382.4
$ws.Cells.Item(64, 12).Value = 3505147.2
$ws.Cells.Item(64, 14).Value = -3505687.2
$ws.Cells.Item(67, 8).Value = 1001756.3
$ws.Cells.Item(67, 10).Value = 1168382.4
$ws.Cells.Item(67, 12).Value = 3505147.2
$ws.Cells.Item(67, 14).Value = -3507019.2
$ws.Cells.Item(69, 8).Value = 5000
$ws.Cells.Item(69, 9).Value = 5000
$ws.Cells.Item(69, 10).Value = 5000
$ws.Cells.Item(69, 11).Value = 15000
$ws.Cells.Item(69, 12).Value = 15000
$ws.Cells.Item(69, 13).Value = -14189
$ws.Cells.Item(69, 14).Value = -16622
$ws.Cells.Item(70, 8).Value = 1799.7333
$ws.Cells.Item(70, 10).Value = 1799.7333
$ws.Cells.Item(70, 12).Value = 5399.199900000001
$ws.Cells.Item(70, 14).Value = -6029.199900000001
$ws.Cells.Item(72, 8).Value = 5000
$ws.Cells.Item(72, 9).Value = 5000
$ws.Cells.Item(72, 10).Value = 5000
$ws.Cells.Item(72, 11).Value = 45000
$ws.Cells.Item(72, 12).Value = 45000
$ws.Cells.Item(72, 13).Value = -40944
$ws.Cells.Item(72, 14).Value = -53112
$ws.Cells.Item(73, 8).Value = 1799.7333
$ws.Cells.Item(73, 10).Value = 1799.7333
$ws.Cells.Item(73, 12).Value = 5399.199900000001
$ws.Cells.Item(73, 14).Value = -7583.199900000001
$ws.Cells.Item(76, 8).Value = 18499.666
$ws.Cells.Item(76, 9).Value = 10998
$ws.Cells.Item(76, 11).Value = 32994
$ws.Cells.Item(76, 13).Value = -32611
$ws.Cells.Item(79, 8).Value = 18499.666
$ws.Cells.Item(79, 9).Value = 10998
$ws.Cells.Item(79, 11).Value = 32994
$ws.Cells.Item(79, 13).Value = -31668
$ws.Cells.Item(87, 8).Value = 14624.25
$ws.Cells.Item(87, 9).Value = 14338
$ws.Cells.Item(87, 11).Value = 43014
$ws.Cells.Item(87, 13).Value = -41766
$ws.Cells.Item(90, 8).Value = 14624.25
$ws.Cells.Item(90, 9).Value = 14338
$ws.Cells.Item(90, 11).Value = 129042
$ws.Cells.Item(90, 13).Value = -122802
$ws.Cells.Item(99, 8).Value = 12109.8
$ws.Cells.Item(99, 9).Value = 183
$ws.Cells.Item(99, 11).Value = 549
$ws.Cells.Item(99, 13).Value = 1697
$ws.Cells.Item(138, 8).Value = 38475336
$ws.Cells.Item(138, 9).Value = 111116390
$ws.Cells.Item(138, 11).Value = 333349170
$ws.Cells.Item(138, 13).Value = -333344030
$ws.Cells.Item(139, 8).Value = 8339626
$ws.Cells.Item(139, 9).Value = 16670262
$ws.Cells.Item(139, 10).Value = 8989.9
$ws.Cells.Item(139, 11).Value = 50010786
$ws.Cells.Item(139, 12).Value = 26969.7
$ws.Cells.Item(139, 13).Value = -50005646
$ws.Cells.Item(139, 14).Value = -37249.7

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(135, 8).Value = 88666.664
$ws.Cells.Item(135, 10).Value = 88666.664
$ws.Cells.Item(135, 12).Value = 88666.664
$ws.Cells.Item(135, 14).Value = -98806.664

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2043.6
$ws.Cells.Item(61, 9).Value = 2178.6667
$ws.Cells.Item(61, 11).Value = 2178.6667
$ws.Cells.Item(61, 13).Value = -1976.6667
$ws.Cells.Item(82, 8).Value = 1487.75
$ws.Cells.Item(82, 9).Value = 913.75
$ws.Cells.Item(82, 10).Value = 1774.75
$ws.Cells.Item(82, 11).Value = 913.75
$ws.Cells.Item(82, 12).Value = 1774.75
$ws.Cells.Item(82, 13).Value = -552.75
$ws.Cells.Item(82, 14).Value = -2496.75
$ws.Cells.Item(85, 8).Value = 1487.75
$ws.Cells.Item(85, 9).Value = 913.75
$ws.Cells.Item(85, 10).Value = 1774.75
$ws.Cells.Item(85, 11).Value = 913.75
$ws.Cells.Item(85, 12).Value = 1774.75
$ws.Cells.Item(85, 13).Value = 334.25
$ws.Cells.Item(85, 14).Value = -4270.75
$ws.Cells.Item(93, 8).Value = 2086.3572
$ws.Cells.Item(93, 9).Value = 1941.4
$ws.Cells.Item(93, 11).Value = 1941.4
$ws.Cells.Item(93, 13).Value = -693.4000000000001
$ws.Cells.Item(113, 8).Value = 2043.6
$ws.Cells.Item(113, 9).Value = 2178.6667
$ws.Cells.Item(113, 11).Value = 2178.6667
$ws.Cells.Item(113, 13).Value = -8.666700000000219

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 48827.168
$ws.Cells.Item(70, 10).Value = 49747.25
$ws.Cells.Item(70, 12).Value = 49747.25
$ws.Cells.Item(70, 14).Value = -50377.25
$ws.Cells.Item(73, 8).Value = 48827.168
$ws.Cells.Item(73, 10).Value = 49747.25
$ws.Cells.Item(73, 12).Value = 49747.25
$ws.Cells.Item(73, 14).Value = -51931.25
$ws.Cells.Item(132, 8).Value = 1699.7894
$ws.Cells.Item(132, 9).Value = 1768.5333
$ws.Cells.Item(132, 10).Value = 1442
$ws.Cells.Item(132, 11).Value = 5305.5999
$ws.Cells.Item(132, 12).Value = 4326
$ws.Cells.Item(132, 13).Value = -2775.5999
$ws.Cells.Item(132, 14).Value = -9386
